$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text ("@") number format on the Price/Volume columns so that
# values like "593.24" or "0.149" are written back as text (matching
# the original inlineStr cell type) instead of being auto-converted
# to numbers by Excel's type inference.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '63.000.08'
$ws.Range('E2').Value = '  +3.08%  '
$ws.Range('D3').Value = '2.953.29'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '593.24'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '148.34'
$ws.Range('E6').Value = '  +3.25%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '2.951.01'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('E9').Value = '  +1.54%  '
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('D11').Value = '0.149'
$ws.Range('E11').Value = '  +6.00%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('E13').Value = '  +5.01%  '
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = '3.442.04'
$ws.Range('D17').Value = '62.985.05'
$ws.Range('E17').Value = '  +3.11%  '
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('D19').Value = '2.952.18'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('D20').Value = '440.81'
$ws.Range('E20').Value = '  +2.11%  '
$ws.Range('D21').Value = '13.49'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = '11.27'
$ws.Range('E24').Value = '  +3.94%  '
$ws.Range('D25').Value = '80.86'
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('D27').Value = '11.81'
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('D30').Value = '7.33'
$ws.Range('E30').Value = '  +6.39%  '
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').Value = '0.0000103'
$ws.Range('E32').Value = '  +18.01%  '
$ws.Range('D33').Value = '26.40'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').Value = '5.63'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').Value = '3.08'
$ws.Range('E38').Value = '  +4.43%  '
$ws.Range('D39').Value = '49.78'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('E40').Value = '  +2.55%  '
$ws.Range('D41').Value = '8.50'
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('E42').Value = '  -3.38%  '
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('D44').Value = '39.52'
$ws.Range('E44').Value = '  -6.17%  '
$ws.Range('D45').Value = '2.709.76'
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('E46').Value = '  +1.59%  '
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('D48').Value = '359.84'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '22.75'
$ws.Range('E51').Value = '  -2.92%  '

# Restore the default (Normal) style so no residual text-format
# styling is left behind on cells that didn't have one originally.
$ws.Range('D2:E51').Style = 'Normal'
